# Insert a new data row at row 130 (pushes existing rows 130..250 down to
# 131..251, extending the used range to A1:R251), then populate the new
# row with the new "Ajo" (garlic) price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(130).Insert()

$ws.Cells.Item(130, 1).Value = 3
$ws.Cells.Item(130, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(130, 3).Value = 'Coquimbo'
$ws.Cells.Item(130, 4).Value = 44512
$ws.Cells.Item(130, 5).Value = 5
$ws.Cells.Item(130, 6).Value = 100112003
$ws.Cells.Item(130, 7).Value = 'Ajo'
$ws.Cells.Item(130, 8).Value = 'Chino'
$ws.Cells.Item(130, 9).Value = 'Primera'
$ws.Cells.Item(130, 10).Value = 65
$ws.Cells.Item(130, 11).Value = 16000
$ws.Cells.Item(130, 12).Value = 16500
$ws.Cells.Item(130, 13).Value = 16231
$ws.Cells.Item(130, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(130, 15).Value = 'China'
$ws.Cells.Item(130, 16).Value = 1623
$ws.Cells.Item(130, 17).Value = 10
$ws.Cells.Item(130, 18).Value = 'Hortaliza'
